# The deck's theme (ppt/theme/theme2.xml, used by the slide master /
# presentation) currently carries the "Integral" color scheme. The edit
# swaps it for the default "Office Theme" 12-slot color scheme (the
# palette that used to live on the Notes Master's theme part):
#
#   dk1/lt1       unchanged (black / white)
#   dk2           455F51 -> 44546A
#   lt2           E3DED1 -> E7E6E6
#   accent1       99CB38 -> 5B9BD5
#   accent2       63A537 -> ED7D31
#   accent3       E6D024 -> A5A5A5
#   accent4       CC9700 -> FFC000
#   accent5       4EB3CF -> 4472C4
#   accent6       378DA6 -> 70AD47
#   hlink         6B9F25 -> 0563C1
#   folHlink      B26B02 -> 954F72
#
# Go through the first slide's ThemeColorScheme (a live view onto the
# presentation's shared theme - every slide/master shares the same
# twelve slots) and set each of the twelve scheme slots via its RGB
# property, exactly like recording "Design > Variants > Colors >
# Customize Colors" in the PowerPoint UI would produce.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-SchemeColor($themeColors, $index, $r, $g, $b) {
    $themeColors.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-SchemeColor $tcs 1  0x00 0x00 0x00   # dk1
Set-SchemeColor $tcs 2  0xFF 0xFF 0xFF   # lt1
Set-SchemeColor $tcs 3  0x44 0x54 0x6A   # dk2
Set-SchemeColor $tcs 4  0xE7 0xE6 0xE6   # lt2
Set-SchemeColor $tcs 5  0x5B 0x9B 0xD5   # accent1
Set-SchemeColor $tcs 6  0xED 0x7D 0x31   # accent2
Set-SchemeColor $tcs 7  0xA5 0xA5 0xA5   # accent3
Set-SchemeColor $tcs 8  0xFF 0xC0 0x00   # accent4
Set-SchemeColor $tcs 9  0x44 0x72 0xC4   # accent5
Set-SchemeColor $tcs 10 0x70 0xAD 0x47   # accent6
Set-SchemeColor $tcs 11 0x05 0x63 0xC1   # hlink
Set-SchemeColor $tcs 12 0x95 0x4F 0x72   # folHlink
